# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a stale "Strike#" style value; this
# recalculated/regenerated set of strikeout (K) counts for each outing and
# rewrites column G (rows 2-21) with the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 4
    4  = 5
    5  = 2
    6  = 3
    7  = 2
    8  = 1
    9  = 4
    10 = 4
    11 = 3
    12 = 0
    13 = 1
    14 = 5
    15 = 3
    16 = 2
    17 = 4
    18 = 6
    19 = 4
    20 = 3
    21 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
